# Apply the two changes captured by the target diff:
#
# 1) The table on slide 16 (the "C1/C2/C3" table) switches from the
#    deck's custom "Table_0" style to a different (built-in) table
#    style id.
#
# 2) The deck's theme color scheme is swapped from "Integral" to the
#    stock "Office Theme" palette (the font scheme / format scheme were
#    already identical between the two theme parts, so only the color
#    values actually change).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 16 -------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6C96DBC3-1738-4E95-9A94-A3850FDACD7E}")
    }
}

# --- 2) Theme color scheme: Integral -> Office Theme -----------------
# Slots are in <a:clrScheme> document order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $rgbHex = $officeThemeColors[$i - 1]
    $r = [math]::Floor($rgbHex / 65536) % 256
    $g = [math]::Floor($rgbHex / 256) % 256
    $b = $rgbHex % 256
    $oleColor = $b * 65536 + $g * 256 + $r
    $themeColors.Colors($i).RGB = $oleColor
}
